# Update BOM rows for the changed Mosfet (Q1, Q5) and sense resistor
# (R7, R10, R15, R17, R18, R22, R25, R26) parts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Einsy Rambo_projectVersion")

# --- Row 34: Q1, Q5 Mosfet -> replaced with new part (PSMN1R8-40YLC) ---
$ws.Cells.Item(34, 2).Value  = "Digikey"
$ws.Cells.Item(34, 3).Value  = "1727-1052-2-ND"
$ws.Cells.Item(34, 5).Value  = "PSMN1R8-40YLC,115"
$ws.Cells.Item(34, 6).Value  = ""
$ws.Cells.Item(34, 9).Value  = "PSMN1R8-40YLC"
$ws.Cells.Item(34, 11).Value = "MOSFET N-CH 40V 100A LFPAK"
$ws.Cells.Item(34, 12).Value = "N-Channel 40V 100A (Tc) 272W (Tc) Surface Mount LFPAK56, Power-SO8"
$ws.Cells.Item(34, 13).Value = "N-FET"
$ws.Cells.Item(34, 14).Value = "-NA-"
$ws.Cells.Item(34, 15).Value = "40V 100A 1.8mOhm"

# --- Row 41: R7, R10, R15, R17, R18, R22, R25, R26 sense resistor -> new part (RUT3216FR220CS) ---
$ws.Cells.Item(41, 3).Value  = "RUT3216FR220CS"
$ws.Cells.Item(41, 5).Value  = "RUT3216FR220CS"
$ws.Cells.Item(41, 11).Value = "RES SMD 0.22 OHM 1% 1/3W 1206"
$ws.Cells.Item(41, 12).Value = "RES SMD 0.22 OHM 1% 1/3W 1206"
